$d = $word.ActiveDocument

# 1) Professional Summary paragraph text replacement
$old1 = "Accomplished Data Scientist with a Bachelor of Science in Computer Science from Stanford University, specializing in machine learning, predictive analytics, and data visualization. With a proven track record of leveraging advanced analytical techniques to enhance business decision-making and insights, I am adept at implementing complex models and interpreting large datasets. My expertise includes proficiency in Python, SQL, R, AWS, and Big Data Technologies. I excel in collaborative environments, working effectively with cross-functional teams to transform business challenges into data-driven solutions. I am eager to bring my strong analytical skills and innovative approach to a dynamic team to drive business success."
$new1 = "Accomplished Data Scientist with a Bachelor of Science in Computer Science from Stanford University, specializing in machine learning, predictive analytics, and data visualization. I have a proven track record in utilizing advanced analytics to drive strategic business decisions and enhance operational efficiency. My expertise includes building sophisticated data models, developing actionable insights through robust reporting and dashboarding, and leveraging tools like Python, SQL, AWS, and Tableau. I am adept at collaborating with cross-functional teams to deliver high-quality, data-driven solutions that align with business goals. My professional experience and strong analytical skills make me a valuable asset for achieving exceptional business insights and improvements."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2) Skills paragraph text replacement
$old2 = "Python, Machine Learning, Predictive Analytics, SQL, Data Visualization, AWS, Statistical Analysis, R, Big Data Technologies, Tableau"
$new2 = "Python, SQL, AWS, Tableau, Machine Learning, Predictive Analytics, Data Visualization, Statistical Analysis"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) Remove the "Junior Data Scientist at Insight Data Science" heading + its 3 bullet points
$old3 = "Junior Data Scientist at Insight Data Science (Mar 2023-Present)"
$headingPar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq $old3) {
        $headingPar = $i
        break
    }
}
if ($headingPar -ne $null) {
    $startRange = $d.Paragraphs.Item($headingPar).Range.Start
    $endRange = $d.Paragraphs.Item($headingPar + 3).Range.End
    $d.Range($startRange, $endRange).Delete() | Out-Null
}

# 4) Update the three "Data Analyst at Tech Solutions" bullet points
$old4 = "Designed and implemented interactive dashboards using Tableau, improving the accessibility and understanding of business metrics across departments."
$new4 = "Developed and maintained advanced analytics dashboards and KPIs, significantly enhancing the monitoring and optimization of business processes."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

$old5 = "Employed advanced statistical techniques and machine learning algorithms to forecast sales trends, which helped increase accuracy in inventory management by 25%."
$new5 = "Utilized machine learning techniques to refine predictive analytics capabilities, resulting in a marked improvement in business insights and decision-making."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

$old6 = "Collaborated with the marketing team to utilize predictive analytics in targeting potential customers, boosting marketing campaign effectiveness by 30%."
$new6 = "Led cross-functional teams in the implementation of strategic data-driven initiatives, contributing to substantial improvements in operational efficiencies and profitability."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# 5) Append the new "Intern Data Scientist at Innovative Startups Inc." section at the end
$lastPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPar.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Heading2"
$p.Range.Text = "Intern Data Scientist at Innovative Startups Inc. (June 2020-Dec 2020)"

$lastPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPar.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "ListBullet"
$p.Range.Text = "Played a key role in the development and enhancement of machine learning models for big data analysis, boosting predictive accuracy by 20%."

$lastPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPar.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "ListBullet"
$p.Range.Text = "Collaborated with senior data scientists to streamline model performance, significantly speeding up data processing and analysis."

$lastPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPar.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "ListBullet"
$p.Range.Text = "Developed customized data visualization tools that improved the communication of complex analytical results to non-technical stakeholders, enhancing understanding and engagement."

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
